# Adding search test cases
#
# "Test Cases" sheet: two existing rows (96/97) pick up a Results column
# (previously absent) and their styling is normalised to match the rest of
# the table; two brand-new rows (98/99) are appended for:
#   TestCase_B97 / OPQA-565 - no filtering options on ALL search results
#   TestCase_B98 / OPQA-571 - search drop down content type retained

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$xlPasteFormats = -4122

# Reference cells already carrying the formatting we need elsewhere in the
# sheet, used purely as format-paint sources (values are left untouched).
$fmtTcid        = $ws.Cells.Item(8, 1)    # A8  - bordered, no fill, no wrap
$fmtJira        = $ws.Cells.Item(71, 2)   # B71 - bordered, no fill, wrap
$fmtDescription = $ws.Cells.Item(90, 3)   # C90 - bordered, no fill, wrap
$fmtRunmode     = $ws.Cells.Item(8, 1)    # A8  - bordered, no fill, no wrap
$fmtResults     = $ws.Cells.Item(93, 5)   # E93 - bordered, no fill, no wrap

# --- Normalise rows 96 & 97 (Description / Runmode / new Results cell) ---
foreach ($r in 96, 97) {
    $fmtDescription.Copy()
    $ws.Cells.Item($r, 3).PasteSpecial($xlPasteFormats)

    $fmtRunmode.Copy()
    $ws.Cells.Item($r, 4).PasteSpecial($xlPasteFormats)

    $fmtResults.Copy()
    $ws.Cells.Item($r, 5).PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = $false

# --- Row 98: TestCase_B97 / OPQA-565 --------------------------------------
$ws.Cells.Item(98, 1).Value = "TestCase_B97"
$ws.Cells.Item(98, 2).Value = "OPQA-565"
$ws.Cells.Item(98, 3).Value = "Verify that no filtering options are present in ALL search results page"
$ws.Cells.Item(98, 4).Value = "Y"
$ws.Cells.Item(98, 5).Value = "PASS"

# --- Row 99: TestCase_B98 / OPQA-571 --------------------------------------
$ws.Cells.Item(99, 1).Value = "TestCase_B98"
$ws.Cells.Item(99, 2).Value = "OPQA-571"
$ws.Cells.Item(99, 3).Value = "Verify that search drop down content type is retained when user navigates back to ALL search results page from record view page"
$ws.Cells.Item(99, 4).Value = "Y"
$ws.Cells.Item(99, 5).Value = "PASS"

# --- Formatting for the two new rows (same look as the rest of the table) -
foreach ($r in 98, 99) {
    $fmtTcid.Copy()
    $ws.Cells.Item($r, 1).PasteSpecial($xlPasteFormats)

    $fmtJira.Copy()
    $ws.Cells.Item($r, 2).PasteSpecial($xlPasteFormats)

    $fmtDescription.Copy()
    $ws.Cells.Item($r, 3).PasteSpecial($xlPasteFormats)

    $fmtRunmode.Copy()
    $ws.Cells.Item($r, 4).PasteSpecial($xlPasteFormats)

    $fmtResults.Copy()
    $ws.Cells.Item($r, 5).PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = $false

# --- View housekeeping: active cell moves with the new last row ----------
$ws.Range("C96").Select()
